$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Some Price cells hold numeric-looking text (e.g. "615.17", "0.0000258",
# "6.00", "0.0400") that Excel would otherwise auto-convert to numbers and
# reformat (losing trailing zeros / switching to scientific notation), so
# those specific cells are forced to Text format before assignment.

$ws.Range("D2").Value = "64.580.31"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "3.159.14"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.17"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.73"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D8").Value = "3.155.97"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").Value = "3.678.83"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "64.612.61"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "3.161.93"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.93"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.14"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.02"
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.75"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.40"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.55"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  -8.44%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.51"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "0.0₃0778"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.00"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "53.04"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "460.52"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0400"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").Value = "  -4.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.41"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "2.845.97"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.269"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.60"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.34"
$ws.Range("E51").Value = "  +0.43%  "
